$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "304.74"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.24%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.87"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.28%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.065"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.16%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08087"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.32%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.920"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.32%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.151"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.51%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "7.839"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.16%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9316"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.44%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1255"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-13.01%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1920"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.13%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09233"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2.25%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03487"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.39%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09896"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.45%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001423"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.85%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006683"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "14.92%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.617"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2.32%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.233"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "8.06%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3423"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.07%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1336"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.168"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.63%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04403"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.85%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001234"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.76%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004720"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.97%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001300"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "5.67%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003127"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "3.45%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01978"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "8.06%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05170"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "8.78%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007530"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.74%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01011"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-3.96%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1365"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "2.84%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002099"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.50%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01067"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-2.10%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006329"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.71%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.02%"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.86%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001599"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-3.61%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.02%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.02%"
